$d = $word.ActiveDocument

# 1. Delete the first paragraph (Harrison Schramm bio), including its paragraph mark.
$d.Paragraphs(1).Range.Delete()

# 2. Remove "Look here to see the example working." (the hyperlinked "here") from the
#    Conclusion paragraph, leaving "...graphics and analysis. " intact.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Look here to see the example working.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
